# Updates the "Estado de Cuenta" worker/period table (rows 16-29) so that
# the two workers (URBIS AMAYA ACONCHA / 8643805 and LUIS FERNANDO VANEGAS
# BASANTA / 1042427825) are interleaved row-by-row and their mora periods
# run in ascending order (2010, 2011, 2012, 2101, 2102, 2103, 2104), per
# the "Actualiza base de datos EC y agrega parte 1 de nuevos estado de
# cuenta" update.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$worker1Doc    = "8643805"
$worker1Name   = "URBIS AMAYA ACONCHA"
$worker1Fee    = 68448
$worker1Base   = 1711200

$worker2Doc    = "1042427825"
$worker2Name   = "LUIS FERNANDO VANEGAS BASANTA"
$worker2Fee    = 40344
$worker2Base   = 1008600

$periods = @("2010", "2011", "2012", "2101", "2102", "2103", "2104")

$row = 16
foreach ($period in $periods) {

    $isLast = ($period -eq "2104")

    # Worker 1 row
    $ws.Cells.Item($row, 2).Value = "CC"
    $ws.Cells.Item($row, 3).Value = $worker1Doc
    $ws.Cells.Item($row, 4).Value = $worker1Name
    $ws.Cells.Item($row, 5).Value = $period
    if ($isLast) {
        $ws.Cells.Item($row, 6).Value = 59322
    } else {
        $ws.Cells.Item($row, 6).Value = $worker1Fee
    }
    $ws.Cells.Item($row, 7).Value = $worker1Base
    $row = $row + 1

    # Worker 2 row
    $ws.Cells.Item($row, 2).Value = "CC"
    $ws.Cells.Item($row, 3).Value = $worker2Doc
    $ws.Cells.Item($row, 4).Value = $worker2Name
    $ws.Cells.Item($row, 5).Value = $period
    if ($isLast) {
        $ws.Cells.Item($row, 6).Value = 34965
    } else {
        $ws.Cells.Item($row, 6).Value = $worker2Fee
    }
    $ws.Cells.Item($row, 7).Value = $worker2Base
    $row = $row + 1
}
